# Updates cryptos list values to match latest snapshot (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force the cell to remain plain text even when the value looks numeric
    # (e.g. "1.003"), then restore the default "Normal" style so no stray
    # number-format style lingers on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '25.978.74'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '1.756.90'
$ws.Range("E3").Value = '  +0.22%  '

Set-TextCell 'D4' '1.003'
$ws.Range("E4").Value = '  +0.24%  '

Set-TextCell 'D5' '236.83'
$ws.Range("E5").Value = '  -0.89%  '

Set-TextCell 'D6' '1.003'
$ws.Range("E6").Value = '  +0.43%  '

Set-TextCell 'D7' '0.5207'
$ws.Range("E7").Value = '  +2.17%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D8' '40.62'
$ws.Range("E8").Value = '  -4.08%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 'D9' '0.2733'
$ws.Range("E9").Value = '  -1.96%  '

Set-TextCell 'D10' '0.06169'
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").Value = '1.772.40'
$ws.Range("E11").Value = '  +1.16%  '

Set-TextCell 'D12' '0.07029'
$ws.Range("E12").Value = '  +0.94%  '

Set-TextCell 'D13' '15.57'
$ws.Range("E13").Value = '  -1.51%  '

Set-TextCell 'D14' '0.6360'
$ws.Range("E14").Value = '  +4.08%  '

Set-TextCell 'D15' '4.520'
$ws.Range("E15").Value = '  -0.14%  '

Set-TextCell 'D16' '77.74'
$ws.Range("E16").Value = '  +0.07%  '

Set-TextCell 'D17' '1.002'
$ws.Range("E17").Value = '  +0.21%  '

Set-TextCell 'D18' '1.004'
$ws.Range("E18").Value = '  +0.43%  '

$ws.Range("D19").Value = '25.984.02'
$ws.Range("E19").Value = '  +0.43%  '

Set-TextCell 'D20' '11.59'
$ws.Range("E20").Value = '  -0.80%  '

Set-TextCell 'D21' '0.000006695'
$ws.Range("E21").Value = '  -3.88%  '

$ws.Range("D22").Value = '2.000.24'
$ws.Range("E22").Value = '  +1.56%  '

Set-TextCell 'D23' '4.060'
$ws.Range("E23").Value = '  -0.63%  '

Set-TextCell 'D24' '8.439'
$ws.Range("E24").Value = '  +2.56%  '

Set-TextCell 'D25' '5.167'
$ws.Range("E25").Value = '  -1.92%  '

$ws.Range("E26").Value = '  +0.75%  '

Set-TextCell 'D27' '1.509'
$ws.Range("E27").Value = '  +1.90%  '

Set-TextCell 'D28' '1.837'
$ws.Range("E28").Value = '  +0.57%  '

Set-TextCell 'D29' '15.15'
$ws.Range("E29").Value = '  +0.68%  '

Set-TextCell 'D30' '103.11'
$ws.Range("E30").Value = '  -0.66%  '

Set-TextCell 'D31' '0.08369'
$ws.Range("E31").Value = '  +1.87%  '

Set-TextCell 'D32' '3.666'
$ws.Range("E32").Value = '  -1.20%  '

Set-TextCell 'D33' '3.428'
$ws.Range("E33").Value = '  -2.54%  '

Set-TextCell 'D34' '0.04463'
$ws.Range("E34").Value = '  -1.44%  '

Set-TextCell 'D35' '2.629'
$ws.Range("E35").Value = '  -0.20%  '

Set-TextCell 'D36' '0.9952'
$ws.Range("E36").Value = '  +0.22%  '

Set-TextCell 'D37' '0.6043'
$ws.Range("E37").Value = '  -1.65%  '

Set-TextCell 'D38' '2.712'
$ws.Range("E38").Value = '  +0.57%  '

Set-TextCell 'D39' '0.01591'
$ws.Range("E39").Value = '  +2.21%  '

Set-TextCell 'D40' '1.960'
$ws.Range("E40").Value = '  +3.01%  '

Set-TextCell 'D41' '1.003'
$ws.Range("E41").Value = '  +0.45%  '

Set-TextCell 'D42' '102.58'
$ws.Range("E42").Value = '  -1.01%  '

Set-TextCell 'D43' '0.3860'
$ws.Range("E43").Value = '  -0.57%  '

Set-TextCell 'D44' '0.7434'
$ws.Range("E44").Value = '  -0.46%  '

Set-TextCell 'D45' '4.925'
$ws.Range("E45").Value = '  -0.18%  '

Set-TextCell 'D46' '0.05507'
$ws.Range("E46").Value = '  +1.81%  '

Set-TextCell 'D47' '6.305'
$ws.Range("E47").Value = '  +4.85%  '

Set-TextCell 'D48' '0.1114'
$ws.Range("E48").Value = '  -0.14%  '

Set-TextCell 'D49' '30.10'
$ws.Range("E49").Value = '  -0.46%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D50' '52.41'
$ws.Range("E50").Value = '  -0.89%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextCell 'D51' '1.006'
$ws.Range("E51").Value = '  +1.07%  '
